$d = $word.ActiveDocument

# Locate the paragraph that ends the "An assumption was made..." sentence; the
# three new paragraphs are inserted right after it (and before "Design Decisions").
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*It was not clear in the task if this should be the case.*") {
        $anchor = $p
        break
    }
}

$r = $anchor.Range

# Insert 5 blank paragraphs after the anchor:
#   base+0, base+1, base+2 -> will hold the 3 runs of the "An issue ... comments." paragraph
#   base+3                 -> "An issue may only have one assigned user."
#   base+4                 -> "StartDate and EndDate filters for GetIssues only cares about Date and not exact time."
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

$base = $anchor.Index + 1

$d.Paragraphs($base).Range.Text = "An issue "
$d.Paragraphs($base + 1).Range.Text = "may"
$d.Paragraphs($base + 2).Range.Text = " have several comments."
$d.Paragraphs($base + 3).Range.Text = "An issue may only have one assigned user."
$d.Paragraphs($base + 4).Range.Text = "StartDate and EndDate filters for GetIssues only cares about Date and not exact time."

# Join the first three paragraphs into a single paragraph made of three runs by
# deleting the paragraph mark between each of them (this preserves the runs as
# distinct <w:r> elements, matching how the text was actually typed/edited).
$end1 = $d.Paragraphs($base).Range.End
$d.Range($end1 - 1, $end1).Delete()

$end2 = $d.Paragraphs($base).Range.End
$d.Range($end2 - 1, $end2).Delete()

Write-Output ("Paragraph 1: [" + $d.Paragraphs($base).Range.Text + "]")
Write-Output ("Paragraph 2: [" + $d.Paragraphs($base + 1).Range.Text + "]")
Write-Output ("Paragraph 3: [" + $d.Paragraphs($base + 2).Range.Text + "]")
